$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'track tights'
$ws.Range("A2").Value = 'wintergear compression men'
$ws.Range("A3").Value = 'basketball training tights'
$ws.Range("A4").Value = 'basketball protector'
$ws.Range("A5").Value = 'winter leggings for men'
$ws.Range("A6").Value = 'kids tights with knee pads'
$ws.Range("A7").Value = 'ropa de monta?a hombre'
$ws.Range("A8").Value = 'men workout leggings nike'
$ws.Range("A9").Value = 'boys white knee pads basketball'
$ws.Range("A10").Value = 'coyote brown pants with knee pads'
$ws.Range("A11").Value = 'black basketball knee pads nike'
$ws.Range("A12").Value = 'basketball knee pads womens'
$ws.Range("A13").Value = 'cold weather workout pants for men'
$ws.Range("A14").Value = 'mens under armour long underwear pants'
$ws.Range("A15").Value = 'calf compression leggings men'
$ws.Range("A16").Value = 'tights for men nike'
$ws.Range("A17").Value = 'swimming pants for men'
$ws.Range("A18").Value = 'copper compression pants for men'
$ws.Range("A19").Value = 'tesla thermal pants'
$ws.Range("A20").Value = 'hex knee sleeve'
$ws.Range("A21").Value = 'thermal nike'
$ws.Range("A22").Value = 'running base layer men'
$ws.Range("A23").Value = 'long spandex men'
$ws.Range("A24").Value = 'cold gear compression pants men'
$ws.Range("A25").Value = 'soccer winter gear'
$ws.Range("A26").Value = 'kids soccor gear'
$ws.Range("A27").Value = 'mens pants with knee pads'
$ws.Range("A28").Value = 'mens thermal baselayer'
$ws.Range("A29").Value = 'knee pad pants men'
$ws.Range("A30").Value = 'adidas compression pants men'
$ws.Range("A31").Value = 'adult clothing protector'
$ws.Range("A32").Value = 'adult football knee pads'
$ws.Range("A33").Value = 'adult football pants with pads'
$ws.Range("A34").Value = 'athletic knee pads'
$ws.Range("A35").Value = 'athletic tights men'
$ws.Range("A36").Value = 'baseball knee pad'
$ws.Range("A37").Value = 'baseball knee pads'
$ws.Range("A38").Value = 'baseball pants adult small'
$ws.Range("A39").Value = 'baskerball tights'
$ws.Range("A40").Value = 'basketball compression knee pads'
$ws.Range("A41").Value = 'basketball compression pants youth with knee pads'
$ws.Range("A42").Value = 'basketball for youth'
$ws.Range("A43").Value = 'basketball hip pads'
$ws.Range("A44").Value = 'basketball leggings men'
$ws.Range("A45").Value = 'basketball pants for women'
$ws.Range("A46").Value = 'basketball pants men'
$ws.Range("A47").Value = 'basketball tights boys youth'
$ws.Range("A48").Value = 'basketball tights for men mcdavid'
$ws.Range("A49").Value = 'basketball tights with pads for boys'
$ws.Range("A50").Value = 'best basketball knee pads'
$ws.Range("A51").Value = 'big knee pads'
$ws.Range("A52").Value = 'big man knee pads'
$ws.Range("A53").Value = 'bjj leggings'
$ws.Range("A54").Value = 'black knee pads'
$ws.Range("A55").Value = 'black knee pads for volleyball'
$ws.Range("A56").Value = 'black leggings design'
$ws.Range("A57").Value = 'black mesh leggings capri'
$ws.Range("A58").Value = 'black youth baseball pants'
$ws.Range("A59").Value = 'boys athletic tights basketball'
$ws.Range("A60").Value = 'boys basketball knee pads mcdavid'
$ws.Range("A61").Value = 'boys compression leggings'
$ws.Range("A62").Value = 'boys compression leggings youth'
$ws.Range("A63").Value = 'boys compression tights'
$ws.Range("A64").Value = 'break away basketball pants'
$ws.Range("A65").Value = 'capri leggings medium'
$ws.Range("A66").Value = 'capri mens'
$ws.Range("A67").Value = 'capri tights for men'
$ws.Range("A68").Value = 'capris leggings'
$ws.Range("A69").Value = 'cold gear for football'
$ws.Range("A70").Value = 'cold weather panta'
$ws.Range("A71").Value = 'compression for knee'
$ws.Range("A72").Value = 'compression gear'
$ws.Range("A73").Value = 'compression men pants'
$ws.Range("A74").Value = 'compression pants big and tall men'
$ws.Range("A75").Value = 'compression pants padded knees basketball'
$ws.Range("A76").Value = 'compression shorts 3 4 length men'
$ws.Range("A77").Value = 'compression tights for men'
$ws.Range("A78").Value = 'cycling pants'
$ws.Range("A79").Value = 'dark purple basketball knee pads'
$ws.Range("A80").Value = 'elbow knee pads youth'
$ws.Range("A81").Value = 'excersize equipment for men'
$ws.Range("A82").Value = 'football 3 4 tights'
$ws.Range("A83").Value = 'football knee pads'
$ws.Range("A84").Value = 'football leg pads'
$ws.Range("A85").Value = 'football pants'
$ws.Range("A86").Value = 'football pants adult black'
$ws.Range("A87").Value = 'football pants youth'
$ws.Range("A88").Value = 'g form knee pads youth'
$ws.Range("A89").Value = 'gel knee pads'
$ws.Range("A90").Value = 'gel knee pads for men'
$ws.Range("A91").Value = 'girl knee pads'
$ws.Range("A92").Value = 'girls basketball knee pads'
$ws.Range("A93").Value = 'girls tights with knee pads'
$ws.Range("A94").Value = 'girls volleyball knee pads'
$ws.Range("A96").Value = 'girls youth volleyball knee pads'
$ws.Range("A97").Value = 'gym pants for men'
